# UPDATE data for Washington, D.C.
$wb = $excel.ActiveWorkbook

# Base year 2025: update the raw value for B2 (cost input)
$wsBase = $wb.Worksheets.Item("2025")
$wsBase.Range("B2").Value = 242000.00000000003

# Future years derive their value from the 2025 value with a declining
# discount factor applied based on (1 - 0.1 * step)
$ws2030 = $wb.Worksheets.Item("2030")
$ws2030.Range("B2").Formula = "='2025'!B2*(1-0.1*0.2)"

$ws2035 = $wb.Worksheets.Item("2035")
$ws2035.Range("B2").Formula = "='2025'!B2*(1-0.1*0.4)"

$ws2040 = $wb.Worksheets.Item("2040")
$ws2040.Range("B2").Formula = "='2025'!B2*(1-0.1*0.6)"

$ws2045 = $wb.Worksheets.Item("2045")
$ws2045.Range("B2").Formula = "='2025'!B2*(1-0.1*0.8)"

$ws2050 = $wb.Worksheets.Item("2050")
$ws2050.Range("B2").Formula = "='2025'!B2*(1-0.1*1)"
